$p = $ppt.ActivePresentation
$s4 = $p.Slides.Item(4)
$s4.MoveTo(8)
$moved = $p.Slides.Item(8)
Write-Output ("Moved slide SlideID=" + $moved.SlideID)
$np = $moved.NotesPage
for ($i = 1; $i -le $np.Shapes.Count; $i++) {
    $sh = $np.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        Write-Output ("Shape " + $i + " text: [" + $sh.TextFrame.TextRange.Text + "]")
    }
}
